$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Runtimes")

$ws.Cells.Item(2, 2).Value = 3.8
$ws.Cells.Item(2, 3).Value = 2.4
$ws.Cells.Item(2, 4).Value = 5.5
$ws.Cells.Item(2, 5).Value = 42.5
$ws.Cells.Item(2, 6).Value = 28.15
$ws.Cells.Item(3, 2).Value = 7.45
$ws.Cells.Item(3, 3).Value = 4.95
$ws.Cells.Item(3, 4).Value = 12
$ws.Cells.Item(3, 5).Value = 87.55
$ws.Cells.Item(3, 6).Value = 59.55
$ws.Cells.Item(4, 2).Value = 11.45
$ws.Cells.Item(4, 3).Value = 8.199999999999999
$ws.Cells.Item(4, 4).Value = 19.25
$ws.Cells.Item(4, 5).Value = 137.5
$ws.Cells.Item(4, 6).Value = 92.3
$ws.Cells.Item(5, 2).Value = 15.85
$ws.Cells.Item(5, 3).Value = 10.35
$ws.Cells.Item(5, 4).Value = 25.1
$ws.Cells.Item(5, 5).Value = 180.75
$ws.Cells.Item(5, 6).Value = 126.3
$ws.Cells.Item(6, 2).Value = 19.9
$ws.Cells.Item(6, 3).Value = 14.95
$ws.Cells.Item(6, 4).Value = 30.5
$ws.Cells.Item(6, 5).Value = 221.7
$ws.Cells.Item(6, 6).Value = 162.2
$ws.Cells.Item(7, 2).Value = 23.4
$ws.Cells.Item(7, 3).Value = 18.45
$ws.Cells.Item(7, 4).Value = 40.7
$ws.Cells.Item(7, 5).Value = 280.25
$ws.Cells.Item(7, 6).Value = 196.1
$ws.Cells.Item(8, 2).Value = 27.8
$ws.Cells.Item(8, 3).Value = 20.4
$ws.Cells.Item(8, 4).Value = 45.6
$ws.Cells.Item(8, 5).Value = 326.7
$ws.Cells.Item(8, 6).Value = 231.65
$ws.Cells.Item(9, 2).Value = 33
$ws.Cells.Item(9, 3).Value = 23.05
$ws.Cells.Item(9, 4).Value = 52.6
$ws.Cells.Item(9, 5).Value = 367.9
$ws.Cells.Item(9, 6).Value = 264.9
$ws.Cells.Item(10, 2).Value = 35.55
$ws.Cells.Item(10, 3).Value = 25.25
$ws.Cells.Item(10, 4).Value = 57.9
$ws.Cells.Item(10, 5).Value = 410.25
$ws.Cells.Item(10, 6).Value = 302.65
$ws.Cells.Item(11, 2).Value = 39.3
$ws.Cells.Item(11, 3).Value = 31.5
$ws.Cells.Item(11, 4).Value = 63.75
$ws.Cells.Item(11, 5).Value = 448.9
$ws.Cells.Item(11, 6).Value = 339.75
$ws.Cells.Item(12, 2).Value = 43.4
$ws.Cells.Item(12, 3).Value = 34.55
$ws.Cells.Item(12, 4).Value = 76.05
$ws.Cells.Item(12, 5).Value = 509.25
$ws.Cells.Item(12, 6).Value = 374.7
$ws.Cells.Item(13, 2).Value = 48.3
$ws.Cells.Item(13, 3).Value = 37.15
$ws.Cells.Item(13, 4).Value = 82.59999999999999
$ws.Cells.Item(13, 5).Value = 561.9
$ws.Cells.Item(13, 6).Value = 407.75
$ws.Cells.Item(14, 2).Value = 52.6
$ws.Cells.Item(14, 3).Value = 39
$ws.Cells.Item(14, 4).Value = 90.2
$ws.Cells.Item(14, 5).Value = 614.1
$ws.Cells.Item(14, 6).Value = 445.25
$ws.Cells.Item(15, 2).Value = 57.45
$ws.Cells.Item(15, 3).Value = 42.2
$ws.Cells.Item(15, 4).Value = 95.09999999999999
$ws.Cells.Item(15, 5).Value = 660.55
$ws.Cells.Item(15, 6).Value = 482.2
$ws.Cells.Item(16, 2).Value = 63.9
$ws.Cells.Item(16, 3).Value = 44.95
$ws.Cells.Item(16, 4).Value = 100.9
$ws.Cells.Item(16, 5).Value = 705.05
$ws.Cells.Item(16, 6).Value = 520.15
$ws.Cells.Item(17, 2).Value = 66.8
$ws.Cells.Item(17, 3).Value = 47.5
$ws.Cells.Item(17, 4).Value = 109.3
$ws.Cells.Item(17, 5).Value = 749.2
$ws.Cells.Item(17, 6).Value = 558.2
$ws.Cells.Item(18, 2).Value = 70.8
$ws.Cells.Item(18, 3).Value = 51
$ws.Cells.Item(18, 4).Value = 115
$ws.Cells.Item(18, 5).Value = 792.05
$ws.Cells.Item(18, 6).Value = 596.35
$ws.Cells.Item(19, 2).Value = 73.2
$ws.Cells.Item(19, 3).Value = 52.15
$ws.Cells.Item(19, 4).Value = 120.6
$ws.Cells.Item(19, 5).Value = 827.45
$ws.Cells.Item(19, 6).Value = 629.25
$ws.Cells.Item(20, 2).Value = 76.25
$ws.Cells.Item(20, 3).Value = 58.65
$ws.Cells.Item(20, 4).Value = 127.4
$ws.Cells.Item(20, 5).Value = 869.7
$ws.Cells.Item(20, 6).Value = 669.6
$ws.Cells.Item(21, 2).Value = 78.2
$ws.Cells.Item(21, 3).Value = 64.34999999999999
$ws.Cells.Item(21, 4).Value = 134.1
$ws.Cells.Item(21, 5).Value = 913.55
$ws.Cells.Item(21, 6).Value = 708.9
$ws.Cells.Item(22, 2).Value = 80.95
$ws.Cells.Item(22, 3).Value = 68.25
$ws.Cells.Item(22, 4).Value = 140.35
$ws.Cells.Item(22, 5).Value = 960.75
$ws.Cells.Item(22, 6).Value = 742.95
$ws.Cells.Item(23, 2).Value = 88.7
$ws.Cells.Item(23, 3).Value = 70.59999999999999
$ws.Cells.Item(23, 4).Value = 156.05
$ws.Cells.Item(23, 5).Value = 1033.25
$ws.Cells.Item(23, 6).Value = 782.6
$ws.Cells.Item(24, 2).Value = 95.2
$ws.Cells.Item(24, 3).Value = 74.25
$ws.Cells.Item(24, 4).Value = 169.3
$ws.Cells.Item(24, 5).Value = 1108.45
$ws.Cells.Item(24, 6).Value = 826.9
$ws.Cells.Item(25, 2).Value = 100.05
$ws.Cells.Item(25, 3).Value = 78.25
$ws.Cells.Item(25, 4).Value = 178.6
$ws.Cells.Item(25, 5).Value = 1177.75
$ws.Cells.Item(25, 6).Value = 873.2
$ws.Cells.Item(26, 2).Value = 102.5
$ws.Cells.Item(26, 3).Value = 79.05
$ws.Cells.Item(26, 4).Value = 179.4
$ws.Cells.Item(26, 5).Value = 1191.45
$ws.Cells.Item(26, 6).Value = 891.9
$ws.Cells.Item(27, 2).Value = 107.7
$ws.Cells.Item(27, 3).Value = 80.75
$ws.Cells.Item(27, 4).Value = 188.7
$ws.Cells.Item(27, 5).Value = 1250.1
$ws.Cells.Item(27, 6).Value = 934.4
$ws.Cells.Item(28, 2).Value = 114.9
$ws.Cells.Item(28, 3).Value = 83.65000000000001
$ws.Cells.Item(28, 4).Value = 193.4
$ws.Cells.Item(28, 5).Value = 1293.1
$ws.Cells.Item(28, 6).Value = 975.35
$ws.Cells.Item(29, 2).Value = 119.95
$ws.Cells.Item(29, 3).Value = 87.55
$ws.Cells.Item(29, 4).Value = 199.95
$ws.Cells.Item(29, 5).Value = 1348.1
$ws.Cells.Item(29, 6).Value = 1014.2
$ws.Cells.Item(30, 2).Value = 122.7
$ws.Cells.Item(30, 3).Value = 89.34999999999999
$ws.Cells.Item(30, 4).Value = 205.3
$ws.Cells.Item(30, 5).Value = 1389.55
$ws.Cells.Item(30, 6).Value = 1047.3
$ws.Cells.Item(31, 2).Value = 129.1
$ws.Cells.Item(31, 3).Value = 92.45
$ws.Cells.Item(31, 4).Value = 211.6
$ws.Cells.Item(31, 5).Value = 1432.8
$ws.Cells.Item(31, 6).Value = 1085.2
$ws.Cells.Item(32, 2).Value = 133.15
$ws.Cells.Item(32, 3).Value = 96.05
$ws.Cells.Item(32, 4).Value = 219.95
$ws.Cells.Item(32, 5).Value = 1480.65
$ws.Cells.Item(32, 6).Value = 1126.8
$ws.Cells.Item(33, 2).Value = 138.15
$ws.Cells.Item(33, 3).Value = 98.40000000000001
$ws.Cells.Item(33, 4).Value = 225.75
$ws.Cells.Item(33, 5).Value = 1519.5
$ws.Cells.Item(33, 6).Value = 1161.1
$ws.Cells.Item(34, 2).Value = 144
$ws.Cells.Item(34, 3).Value = 101.25
$ws.Cells.Item(34, 4).Value = 234.75
$ws.Cells.Item(34, 5).Value = 1579.35
$ws.Cells.Item(34, 6).Value = 1210.75
$ws.Cells.Item(35, 2).Value = 145.15
$ws.Cells.Item(35, 3).Value = 103.4
$ws.Cells.Item(35, 4).Value = 239.55
$ws.Cells.Item(35, 5).Value = 1610.75
$ws.Cells.Item(35, 6).Value = 1243.8
$ws.Cells.Item(36, 2).Value = 150.35
$ws.Cells.Item(36, 3).Value = 107.6
$ws.Cells.Item(36, 4).Value = 246.85
$ws.Cells.Item(36, 5).Value = 1663.8
$ws.Cells.Item(36, 6).Value = 1289.55
$ws.Cells.Item(37, 2).Value = 151.55
$ws.Cells.Item(37, 3).Value = 109.85
$ws.Cells.Item(37, 4).Value = 254.25
$ws.Cells.Item(37, 5).Value = 1697.3
$ws.Cells.Item(37, 6).Value = 1326.3
$ws.Cells.Item(38, 2).Value = 156.8
$ws.Cells.Item(38, 3).Value = 115.75
$ws.Cells.Item(38, 4).Value = 263.8
$ws.Cells.Item(38, 5).Value = 1762.6
$ws.Cells.Item(38, 6).Value = 1381.4
$ws.Cells.Item(39, 2).Value = 161.6
$ws.Cells.Item(39, 3).Value = 126.75
$ws.Cells.Item(39, 4).Value = 275.85
$ws.Cells.Item(39, 5).Value = 1818.6
$ws.Cells.Item(39, 6).Value = 1434.15
$ws.Cells.Item(40, 2).Value = 164.45
$ws.Cells.Item(40, 3).Value = 135
$ws.Cells.Item(40, 4).Value = 277.9
$ws.Cells.Item(40, 5).Value = 1861.5
$ws.Cells.Item(40, 6).Value = 1477.7
$ws.Cells.Item(41, 2).Value = 167
$ws.Cells.Item(41, 3).Value = 140.05
$ws.Cells.Item(41, 4).Value = 287.25
$ws.Cells.Item(41, 5).Value = 1907.3
$ws.Cells.Item(41, 6).Value = 1514.2
